# Week 1 work started — build out the Portfolio Allocation table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header labels (row 2, columns B:H) ----------------------------------
$ws.Range("B2").Value = "Stock Name"
$ws.Range("C2").Value = "Ticker"
$ws.Range("D2").Value = "Sector"
$ws.Range("E2").Value = "Amount allocated"
$ws.Range("F2").Value = "Price on Buy Date"
$ws.Range("G2").Value = "Shares Bought"
$ws.Range("H2").Value = "Notes"

# ---- Column widths (approximate the original author's best-fit widths) ---
$ws.Columns.Item(2).ColumnWidth = 9.92    # B ~10.78
$ws.Columns.Item(3).ColumnWidth = 10.26   # C ~11.11
$ws.Columns.Item(4).ColumnWidth = 8.75    # D ~9.66
$ws.Columns.Item(5).ColumnWidth = 14.58   # E ~15.55
$ws.Columns.Item(6).ColumnWidth = 14.59   # F ~15.44
$ws.Columns.Item(7).ColumnWidth = 11.58   # G ~12.55

# ---- Row heights -----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15

# ---- Header row formatting (B2:H2): centered, filled, boxed ---------------
$header = $ws.Range("B2:H2")
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4108     # xlCenter
$header.Interior.ThemeColor = 10      # maps to theme index 9 (Accent6) tint ~0.8
$header.Interior.TintAndShade = 0.8

# Outer box (medium) around the whole header band + thin separators between
# each header cell.
$header.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$header.Borders.Item(8).Weight = -4138  # xlMedium
$header.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$header.Borders.Item(9).Weight = -4138  # xlMedium
$header.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$header.Borders.Item(7).Weight = -4138  # xlMedium
$header.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$header.Borders.Item(10).Weight = -4138 # xlMedium
$header.Borders.Item(11).LineStyle = 1  # xlInsideVertical
$header.Borders.Item(11).Weight = 2     # xlThin

# ---- Data rows (3-11): thin vertical rules on B/D/F/H ----------------------
$dataRows = $ws.Range("B3:B11,D3:D11,F3:F11,H3:H11")
$dataRows.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$dataRows.Borders.Item(7).Weight = 2      # xlThin
$dataRows.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$dataRows.Borders.Item(10).Weight = 2     # xlThin

# ---- Last data row (row 12): close the table off with a bottom rule -------
$lastRow = $ws.Range("B12:H12")
$lastRow.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$lastRow.Borders.Item(9).Weight = 2       # xlThin

$lastRowSides = $ws.Range("B12,D12,F12,H12")
$lastRowSides.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$lastRowSides.Borders.Item(7).Weight = 2      # xlThin
$lastRowSides.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$lastRowSides.Borders.Item(10).Weight = 2     # xlThin

# ---- Selection / view state, matching the saved workbook -------------------
$ws.Range("I13").Select()
